# Updated symbol list on Mon Jan 16 10:09:22 UTC 2023 with GitHub Actions
# Refreshes the cryptocurrency price/volume/hour snapshot on Sheet1.
# Each entry below is an A1 cell reference + its new text value. Values are
# written as plain text (matching the sheet's original inline-string
# storage for Price/Volume/Hora) rather than being auto-coerced to
# numbers/percentages by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "298.99" },
    @{ Cell = "E2"; Value = "1.36%" },
    @{ Cell = "G2"; Value = "10" },
    @{ Cell = "D3"; Value = "31.53" },
    @{ Cell = "E3"; Value = "0.88%" },
    @{ Cell = "G3"; Value = "10" },
    @{ Cell = "D4"; Value = "5.160" },
    @{ Cell = "E4"; Value = "1.09%" },
    @{ Cell = "G4"; Value = "10" },
    @{ Cell = "D5"; Value = "0.08026" },
    @{ Cell = "E5"; Value = "8.84%" },
    @{ Cell = "G5"; Value = "10" },
    @{ Cell = "D6"; Value = "2.641" },
    @{ Cell = "E6"; Value = "60.51%" },
    @{ Cell = "G6"; Value = "10" },
    @{ Cell = "D7"; Value = "7.852" },
    @{ Cell = "E7"; Value = "2.19%" },
    @{ Cell = "G7"; Value = "10" },
    @{ Cell = "D8"; Value = "3.827" },
    @{ Cell = "E8"; Value = "2.24%" },
    @{ Cell = "G8"; Value = "10" },
    @{ Cell = "D9"; Value = "0.9078" },
    @{ Cell = "E9"; Value = "-0.80%" },
    @{ Cell = "G9"; Value = "10" },
    @{ Cell = "D10"; Value = "0.1738" },
    @{ Cell = "E10"; Value = "4.23%" },
    @{ Cell = "G10"; Value = "10" },
    @{ Cell = "D11"; Value = "0.07216" },
    @{ Cell = "E11"; Value = "0.57%" },
    @{ Cell = "G11"; Value = "10" },
    @{ Cell = "D12"; Value = "0.08059" },
    @{ Cell = "E12"; Value = "1.03%" },
    @{ Cell = "G12"; Value = "10" },
    @{ Cell = "D13"; Value = "0.03018" },
    @{ Cell = "E13"; Value = "0.94%" },
    @{ Cell = "G13"; Value = "10" },
    @{ Cell = "D14"; Value = "0.09977" },
    @{ Cell = "E14"; Value = "0.79%" },
    @{ Cell = "G14"; Value = "10" },
    @{ Cell = "D15"; Value = "0.001497" },
    @{ Cell = "E15"; Value = "-0.41%" },
    @{ Cell = "G15"; Value = "10" },
    @{ Cell = "D16"; Value = "0.005948" },
    @{ Cell = "E16"; Value = "-4.35%" },
    @{ Cell = "G16"; Value = "10" },
    @{ Cell = "D17"; Value = "3.509" },
    @{ Cell = "E17"; Value = "1.73%" },
    @{ Cell = "G17"; Value = "10" },
    @{ Cell = "D18"; Value = "2.254" },
    @{ Cell = "E18"; Value = "1.18%" },
    @{ Cell = "G18"; Value = "10" },
    @{ Cell = "D19"; Value = "0.3283" },
    @{ Cell = "E19"; Value = "0.19%" },
    @{ Cell = "G19"; Value = "10" },
    @{ Cell = "D20"; Value = "0.1329" },
    @{ Cell = "E20"; Value = "-1.21%" },
    @{ Cell = "G20"; Value = "10" },
    @{ Cell = "D21"; Value = "4.599" },
    @{ Cell = "E21"; Value = "1.11%" },
    @{ Cell = "G21"; Value = "10" },
    @{ Cell = "D22"; Value = "0.1600" },
    @{ Cell = "E22"; Value = "3.33%" },
    @{ Cell = "G22"; Value = "10" },
    @{ Cell = "D23"; Value = "0.04575" },
    @{ Cell = "E23"; Value = "-0.94%" },
    @{ Cell = "G23"; Value = "10" },
    @{ Cell = "D24"; Value = "0.001261" },
    @{ Cell = "E24"; Value = "3.73%" },
    @{ Cell = "G24"; Value = "10" },
    @{ Cell = "D25"; Value = "0.004456" },
    @{ Cell = "E25"; Value = "0.73%" },
    @{ Cell = "G25"; Value = "10" },
    @{ Cell = "D26"; Value = "0.0001180" },
    @{ Cell = "E26"; Value = "-9.07%" },
    @{ Cell = "G26"; Value = "10" },
    @{ Cell = "D27"; Value = "0.0003429" },
    @{ Cell = "E27"; Value = "83.17%" },
    @{ Cell = "G27"; Value = "10" },
    @{ Cell = "G28"; Value = "10" },
    @{ Cell = "G29"; Value = "10" },
    @{ Cell = "G30"; Value = "10" },
    @{ Cell = "G31"; Value = "10" },
    @{ Cell = "G32"; Value = "10" },
    @{ Cell = "G33"; Value = "10" },
    @{ Cell = "G34"; Value = "10" },
    @{ Cell = "G35"; Value = "10" },
    @{ Cell = "G36"; Value = "10" },
    @{ Cell = "G37"; Value = "10" },
    @{ Cell = "G38"; Value = "10" },
    @{ Cell = "D39"; Value = "0.01855" },
    @{ Cell = "E39"; Value = "10.57%" },
    @{ Cell = "G39"; Value = "10" },
    @{ Cell = "D40"; Value = "0.04526" },
    @{ Cell = "E40"; Value = "2.70%" },
    @{ Cell = "G40"; Value = "10" },
    @{ Cell = "D41"; Value = "0.007019" },
    @{ Cell = "E41"; Value = "-0.60%" },
    @{ Cell = "G41"; Value = "10" },
    @{ Cell = "E42"; Value = "1.19%" },
    @{ Cell = "G42"; Value = "10" },
    @{ Cell = "D43"; Value = "0.002240" },
    @{ Cell = "G43"; Value = "10" },
    @{ Cell = "D44"; Value = "0.01043" },
    @{ Cell = "E44"; Value = "-2.42%" },
    @{ Cell = "G44"; Value = "10" },
    @{ Cell = "D45"; Value = "0.00006378" },
    @{ Cell = "E45"; Value = "6.51%" },
    @{ Cell = "G45"; Value = "10" },
    @{ Cell = "D46"; Value = "0.00000000750" },
    @{ Cell = "E46"; Value = "-0.02%" },
    @{ Cell = "G46"; Value = "10" },
    @{ Cell = "B47"; Value = "CoinbaseStockToken" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin" },
    @{ Cell = "D47"; Value = "0.006199" },
    @{ Cell = "E47"; Value = "-43.54%" },
    @{ Cell = "G47"; Value = "10" },
    @{ Cell = "B48"; Value = "BOLO" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo" },
    @{ Cell = "D48"; Value = "0.8206" },
    @{ Cell = "E48"; Value = "-57.44%" },
    @{ Cell = "G48"; Value = "10" },
    @{ Cell = "D49"; Value = "0.00002100" },
    @{ Cell = "E49"; Value = "-0.02%" },
    @{ Cell = "G49"; Value = "10" },
    @{ Cell = "D50"; Value = "0.0002000" },
    @{ Cell = "E50"; Value = "0.05%" },
    @{ Cell = "G50"; Value = "10" },
    @{ Cell = "G51"; Value = "10" }

)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
